$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Price (D) and Volume(1h) (E) columns are text-formatted in the source data,
# (e.g. "63.434.10", "0.999", "  +0.73%  ") so force text number format before
# assigning values that could otherwise be auto-converted to numeric/date types.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.564.01"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.97%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.644.84"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.12%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.29"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.95%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.32"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.85%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.588"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.36%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.644.39"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.12%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.79%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.58"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.28%  "

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.34%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.29%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.44"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.25%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.119.32"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.08%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.441.37"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.87%  "

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.37%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.639.61"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.56%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.36"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.47%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "341.11"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.29%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.36"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.72%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.70"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.04%  "

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.11%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.63"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.58%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.66"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.04%  "

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +5.86%  "

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.67%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "550.74"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +16.78%  "

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.05%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.40"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.71%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.75"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.06%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.82"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +13.22%  "

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.06%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0807"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.08%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "175.36"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.85%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.86"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +7.41%  "

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.89%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.08"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.14%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.77"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.11%  "

# Row 41
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "170.44"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +7.64%  "

# Row 42
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.02%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.25"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.03%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.73"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.79%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.34"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.21%  "

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.04%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0550"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.29%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0959"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.24%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.57%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.65"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.57%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.34"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.72%  "
